$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing City column (column E) to make room
# for the new "Primary Email" column. This shifts City -> F and Fund -> G.
$ws.Columns.Item(5).Insert()

# Match the new column's width to column D's width (19.33203125 in raw units).
$ws.Columns.Item(5).ColumnWidth = 18.5

# Header for the new column
$ws.Range("E1").Value = "Primary Email"
$ws.Range("E1").Style = "Normal 2"

# Email values for each investor row
$ws.Range("E2").Value = "emp1@gmail.com"
$ws.Range("E2").Style = "Normal 2"

$ws.Range("E3").Value = "emp2@gmail.com"
$ws.Range("E3").Style = "Normal 2"

$ws.Range("E4").Value = "emp3@gmail.com"
$ws.Range("E4").Style = "Normal 2"

$ws.Range("E5").Value = "emp4@gmail.com"
$ws.Range("E5").Style = "Normal 2"

$ws.Range("E6").Value = "emp5@gmail.com"
$ws.Range("E6").Style = "Normal 2"

# Row 7 keeps the same "no explicit style" pattern as the other unstyled
# trailing cells in that row (C7/D7), so we leave its style untouched.
$ws.Range("E7").Value = "emp6@gmail.com"

# Normalize column B's style: it was pointing at a duplicate cellXf entry
# (same formatting as "Normal 2"); re-assigning collapses it onto the
# canonical entry and allows the duplicate xf entries to be dropped.
$ws.Range("B1:B7").Style = "Normal 2"

# The two trailing placeholder cells in column B (rows 8-9) are no longer
# needed, clear them completely (value + formatting) so they disappear.
$ws.Cells.Item(8, 2).Clear()
$ws.Cells.Item(9, 2).Clear()

# Update the active selection to reflect the new email column.
$ws.Range("E1:E7").Select()

